$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.929.69"
$ws.Range("E2").Value = "  -4.07%  "
$ws.Range("D3").Value = "1.635.97"
$ws.Range("E3").Value = "  -6.34%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9970"
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.61"
$ws.Range("E5").Value = "  -6.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9995"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4748"
$ws.Range("E7").Value = "  -6.37%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2573"
$ws.Range("E8").Value = "  -6.56%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06084"
$ws.Range("E9").Value = "  -1.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07000"
$ws.Range("E10").Value = "  -3.70%  "
$ws.Range("D11").Value = "1.638.29"
$ws.Range("E11").Value = "  -6.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.58"
$ws.Range("E12").Value = "  -4.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.5975"
$ws.Range("E13").Value = "  -8.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.345"
$ws.Range("E14").Value = "  -6.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "73.46"
$ws.Range("E15").Value = "  -5.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9996"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9974"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").Value = "24.907.57"
$ws.Range("E18").Value = "  -4.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006565"
$ws.Range("E19").Value = "  -4.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.15"
$ws.Range("E20").Value = "  -6.21%  "
$ws.Range("D21").Value = "1.852.46"
$ws.Range("E21").Value = "  -5.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.333"
$ws.Range("E22").Value = "  -2.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.525"
$ws.Range("E23").Value = "  -2.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.218"
$ws.Range("E24").Value = "  -3.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "132.84"
$ws.Range("E25").Value = "  -2.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.80"
$ws.Range("E26").Value = "  -3.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.380"
$ws.Range("E27").Value = "  -8.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "103.33"
$ws.Range("E28").Value = "  -2.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.632"
$ws.Range("E29").Value = "  -8.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.866"
$ws.Range("E30").Value = "  -0.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07687"
$ws.Range("E31").Value = "  -6.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.539"
$ws.Range("E32").Value = "  -3.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9991"
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04303"
$ws.Range("E34").Value = "  -8.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.590"
$ws.Range("E35").Value = "  -2.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9226"
$ws.Range("E36").Value = "  -7.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5785"
$ws.Range("E37").Value = "  -6.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.528"
$ws.Range("E38").Value = "  -8.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8109"
$ws.Range("E41").Value = "  +5.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.28"
$ws.Range("E42").Value = "  -3.73%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.760"
$ws.Range("E43").Value = "  -8.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3681"
$ws.Range("E44").Value = "  -6.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.697"
$ws.Range("E45").Value = "  -6.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.027"
$ws.Range("E48").Value = "  -5.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "29.41"
$ws.Range("E49").Value = "  -4.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.9988"
$ws.Range("E50").Value = "  -0.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9945"
$ws.Range("E51").Value = "  -0.47%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01527"
$ws.Range("E39").Value = "  -5.59%  "
$ws.Range("B40").Value = "PaxDollar"
$ws.Range("C40").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9982"
$ws.Range("E40").Value = "  -0.10%  "
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1089"
$ws.Range("E46").Value = "  -5.78%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05182"
$ws.Range("E47").Value = "  -3.11%  "
